$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rebuild the EC (Estado de Cuenta) detail table rows 16-29.
# Columns: C = N Doc Trabajador, D = Nombre Trabajador, E = Periodo Mora,
#          F = Valor Mora, G = Salario Basico (B = Tipo Doc, stays "CC")

$rows = @(
  @{ Row=16; Doc="1143348923"; Nombre="PABLO EMILIO LONDOÑO RIOS";      Periodo="2309"; Mora=24000; Salario=1160000 },
  @{ Row=17; Doc="1143348923"; Nombre="PABLO EMILIO LONDOÑO RIOS";      Periodo="2304"; Mora=40000; Salario=1160000 },
  @{ Row=18; Doc="1143348923"; Nombre="PABLO EMILIO LONDOÑO RIOS";      Periodo="2303"; Mora=40000; Salario=1160000 },
  @{ Row=19; Doc="1108763580"; Nombre="MARIA CLARA BANQUET ROMERO";     Periodo="2201"; Mora=46400; Salario=1300000 },
  @{ Row=20; Doc="1108763580"; Nombre="MARIA CLARA BANQUET ROMERO";     Periodo="2112"; Mora=46400; Salario=1300000 },
  @{ Row=21; Doc="1108763580"; Nombre="MARIA CLARA BANQUET ROMERO";     Periodo="2111"; Mora=46400; Salario=1300000 },
  @{ Row=22; Doc="1108763580"; Nombre="MARIA CLARA BANQUET ROMERO";     Periodo="2110"; Mora=36341; Salario=908526 },
  @{ Row=23; Doc="1108763580"; Nombre="MARIA CLARA BANQUET ROMERO";     Periodo="2109"; Mora=36341; Salario=908526 },
  @{ Row=24; Doc="1127587489"; Nombre="KAREN CECILIA BARRAGAN MUNZON";  Periodo="2309"; Mora=31200; Salario=1300000 },
  @{ Row=25; Doc="1127587489"; Nombre="KAREN CECILIA BARRAGAN MUNZON";  Periodo="2304"; Mora=46400; Salario=1300000 },
  @{ Row=26; Doc="1127587489"; Nombre="KAREN CECILIA BARRAGAN MUNZON";  Periodo="2303"; Mora=46400; Salario=1300000 },
  @{ Row=27; Doc="1047447875"; Nombre="VICTOR ALFONSO VELASQUEZ IRIARTE"; Periodo="2010"; Mora=35112; Salario=877803 },
  @{ Row=28; Doc="1143391729"; Nombre="EDWIN MANUEL MARTINEZ LOZANO";   Periodo="2304"; Mora=46400; Salario=1423500 },
  @{ Row=29; Doc="1143391729"; Nombre="EDWIN MANUEL MARTINEZ LOZANO";   Periodo="2303"; Mora=46400; Salario=1423500 }
)

foreach ($r in $rows) {
  $ws.Range("C$($r.Row)").Value = $r.Doc
  $ws.Range("D$($r.Row)").Value = $r.Nombre
  $ws.Range("E$($r.Row)").Value = $r.Periodo
  $ws.Range("F$($r.Row)").Value = $r.Mora
  $ws.Range("G$($r.Row)").Value = $r.Salario
}
